$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. SEC_Processes: insert a new process row for the wind turbine
#    (ELE_EX_WIND_TURBINE) right above the existing MIN rows (i.e. new row 8,
#    pushing the old rows 8-13 down by one).
# ---------------------------------------------------------------------------
$secProc = $wb.Worksheets.Item("SEC_Processes")
$secProc.Rows.Item(7).Copy()
$secProc.Rows.Item(8).Insert()
$secProc.Cells.Item(8,4).Value = "ELE_EX_WIND_TURBINE"
$secProc.Cells.Item(8,5).Value = "Wind turbine Onshore"

Write-Host "Step 1 (SEC_Processes) done"
